$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("L21").Value
$ws.Range("Z1").Value = $v
